# Edit: update the "Dt. Referencia" column (G) from 2024-09-09 (45544) to
# 2024-09-10 (45545) for every data row, rename the sheet to reflect the new
# export timestamp, and update a handful of specific values that changed
# between the two export runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the new export run's name/timestamp.
$ws.Name = "IClientBalance-20240910-094026-"

# Column G holds the "Dt. Referencia" serial date for every data row
# (rows 2 through 274). Bump every one of them by one day.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45545
}

# A handful of rows also changed their Saldo Previsto / Vl. Total (and, in
# one case, Vl. Projetado) figures between the two exports.
$ws.Cells.Item(51, 5).Value = 999.99
$ws.Cells.Item(51, 8).Value = 999.99

$ws.Cells.Item(52, 5).Value = 26535.200000000001
$ws.Cells.Item(52, 8).Value = 26535.200000000001

$ws.Cells.Item(104, 4).Value = -13756.32
$ws.Cells.Item(104, 8).Value = 4780.63

$ws.Cells.Item(118, 5).Value = 999.99
$ws.Cells.Item(118, 8).Value = 999.99

$ws.Cells.Item(189, 5).Value = 999.99
$ws.Cells.Item(189, 8).Value = 999.99

$ws.Cells.Item(230, 5).Value = 939.7
$ws.Cells.Item(230, 8).Value = 939.7
